# "Generate Report for Handback"
#
# The handback status report has two per-language sheets ("zh-cn" and
# "de-de"). Each row tracks one file's handoff/handback round-trip; column
# E is "Correspond Handoff Datetime" and column H is "Correspond Handback
# DateTime" (both plain text, formatted like timestamps).
#
# A new handback cycle ran for the first file in each language
# (89f49ef1-...), producing fresh handoff/handback timestamps. The second
# file's row (faa2eff6-...) is unaffected. Update the four timestamp
# cells accordingly.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 16:53:19"
$wsZhCn.Range("H2").Value = "2016-03-22 16:53:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 16:53:23"
$wsDeDe.Range("H2").Value = "2016-03-22 16:53:49"
